$wb = $excel.ActiveWorkbook

# --- Sheet "Metadata": update Date and FHIR Version values ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2025-06-13T15:45:04+00:00"
$meta.Range("B15").Value = "4.0.1"

# --- Sheet "Elements": update constraint text, Type(s) for Extension.id, and R4B -> R4 link ---
$elements = $wb.Worksheets.Item("Elements")

# Extension row (row 2): drop the "unless an empty Parameters resource ... or $this is Parameters" clause
$elements.Range("AJ2").Value = "ele-1:All FHIR elements must have a @value or children {hasValue() or (children().count() > id.count())}`next-1:Must have either extensions or value[x], not both {extension.exists() != value.exists()}"

# Extension.id row (row 3): Type(s) changes from "id" to "string"
$elements.Range("K3").Value = "string`n"

# Extension.value[x] row (row 6): fix FHIR spec link from R4B to R4
$elements.Range("M6").Value = "Value of extension - must be one of a constrained set of the data types (see [Extensibility](http://hl7.org/fhir/R4/extensibility.html) for a list)."
